# Updates question to be different from Q6
# Target slide: the one with the "7kg" weight-force annotation (slide 16).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# 1) "Rectangle 7" becomes a left-pointing arrow, repositioned/resized slightly.
$arrow = $s.Shapes.Item(1)
$arrow.Name = "Left Arrow 7"
$arrow.AutoShapeType = 34            # msoShapeLeftArrow -> prstGeom "leftArrow"
$arrow.Left   = 479.9998788605372    # -> 6095998 EMU
$arrow.Top    = 264.43823283602705   # -> 3358365 EMU
$arrow.Width  = 162.60681300692679   # -> 2065106 EMU
$arrow.Height = 13.550668768377635   # -> 172093 EMU

# 2) Clear the old "7kg" answer text from "Rectangle 9".
$oldAnswer = $s.Shapes.Item(6)
$oldAnswer.TextFrame.TextRange.Text = ""

# 3) Add a new textbox with the new answer, "305N" (with a small leading
#    Wingdings glyph run, matching the original author's formatting).
$newAnswer = $s.Shapes.AddTextbox(1, 1.0, 1.0, 10.0, 10.0)
$tf = $newAnswer.TextFrame
$tr = $tf.TextRange
$tr.Text = " 305N"
$tr.Font.Size = 24
$tr.Font.Bold = $true
$tr.LanguageID = "en-GB"

$glyph = $tr.Characters(1, 1)
$glyph.Font.Name = "Wingdings"

$tf.WordWrap = $false
$tf.AutoSize = 1                     # ppAutoSizeShapeToFitText -> <a:spAutoFit/>

$newAnswer.Left   = 504.07279976516094   # -> 6401724 EMU
$newAnswer.Top    = 273.3299576006947    # -> 3471290 EMU
$newAnswer.Width  = 96.33169239042489    # -> 1223412 EMU
$newAnswer.Height = 36.351613650267396   # -> 461665 EMU
